$wb = $excel.ActiveWorkbook

# --- Sheet R1: remove the LTH0330 row (row 6) ---
$wsR1 = $wb.Worksheets.Item("R1")
$wsR1.Rows.Item(6).Delete()

# --- Sheet R2: bump Elapsed Duration(Hrs) values by 36 seconds ---
$wsR2 = $wb.Worksheets.Item("R2")
$wsR2.Range("G2").Value = "12111:08:55"
$wsR2.Range("G3").Value = "3240:52:24"
$wsR2.Range("G4").Value = "479:03:58"

# --- Sheet R3: add a new outage row (R3 / HAL0947 / SCECO / Weak+Good / Zain) ---
$wsR3 = $wb.Worksheets.Item("R3")
$wsR3.Range("B2").Value = "R3"
$wsR3.Range("D2").Value = "HAL0947"
$wsR3.Range("I2").Value = "SCECO"
$wsR3.Range("J2").Value = "Weak+Good"
$wsR3.Range("L2").Value = "Zain"

# --- Sheet R1: bump Elapsed Duration(Hrs) values by 36 seconds ---
$wsR1.Range("G2").Value = "3929:45:14"
$wsR1.Range("G3").Value = "69:17:52"

# --- Sheet R4: bump Elapsed Duration(Hrs) values by 36 seconds ---
$wsR4 = $wb.Worksheets.Item("R4")
$wsR4.Range("G2").Value = "2956:58:44"
$wsR4.Range("G3").Value = "184:10:59"
$wsR4.Range("G4").Value = "72:23:24"
$wsR4.Range("G5").Value = "70:00:57"

# --- Sheet R5: bump Elapsed Duration(Hrs) value by 36 seconds ---
$wsR5 = $wb.Worksheets.Item("R5")
$wsR5.Range("G2").Value = "430:57:43"

# --- Sheet R6: bump Elapsed Duration(Hrs) value by 36 seconds ---
$wsR6 = $wb.Worksheets.Item("R6")
$wsR6.Range("G2").Value = "71:30:01"
